$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 15627284
$ws.Range("I100").Value = 20002324
$ws.Range("J100").Value = 2142.8572
$ws.Range("K100").Value = 20002324
$ws.Range("L100").Value = 2142.8572
$ws.Range("M100").Value = -20001783
$ws.Range("N100").Value = -3224.8572
$ws.Range("H132").Value = 4766318
$ws.Range("I132").Value = 6175493
$ws.Range("J132").Value = 10351.9375
$ws.Range("K132").Value = 18526479
$ws.Range("L132").Value = 31055.8125
$ws.Range("M132").Value = -18523949
$ws.Range("N132").Value = -36115.8125
$ws.Range("H137").Value = 1094.9595
$ws.Range("I137").Value = 904.7895
$ws.Range("J137").Value = 1295.6945
$ws.Range("K137").Value = 2714.3685
$ws.Range("L137").Value = 3887.0835
$ws.Range("M137").Value = -164.3685
$ws.Range("N137").Value = -8987.083500000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4322.5225
$ws.Range("I32").Value = 3991.3276
$ws.Range("K32").Value = 3991.3276
$ws.Range("M32").Value = -3704.3276
$ws.Range("H45").Value = 1372.0588
$ws.Range("I45").Value = 1346.2142
$ws.Range("J45").Value = 1492.6666
$ws.Range("K45").Value = 1346.2142
$ws.Range("L45").Value = 1492.6666
$ws.Range("M45").Value = -969.2141999999999
$ws.Range("N45").Value = -2246.6666
$ws.Range("H61").Value = 19231718
$ws.Range("I61").Value = 21739926
$ws.Range("K61").Value = 21739926
$ws.Range("M61").Value = -21739714
$ws.Range("H74").Value = 1139.475
$ws.Range("I74").Value = 779.70966
$ws.Range("K74").Value = 779.70966
$ws.Range("M74").Value = 94.29034000000001
$ws.Range("H77").Value = 1139.475
$ws.Range("I77").Value = 779.70966
$ws.Range("K77").Value = 3898.5483
$ws.Range("M77").Value = 469.4517000000001
$ws.Range("H122").Value = 4271
$ws.Range("I122").Value = 3950
$ws.Range("K122").Value = 11850
$ws.Range("M122").Value = -9400
$ws.Range("H132").Value = 3896.125
$ws.Range("I132").Value = 5603
$ws.Range("K132").Value = 16809
$ws.Range("M132").Value = -14279
$ws.Range("H136").Value = 19231718
$ws.Range("I136").Value = 21739926
$ws.Range("K136").Value = 65219778
$ws.Range("M136").Value = -65217228

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 5500
$ws.Range("I37").Value = 1000
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 1000
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = -863
$ws.Range("N37").Value = -10274
$ws.Range("H94").Value = 19232062
$ws.Range("I94").Value = 22728444
$ws.Range("J94").Value = 1955
$ws.Range("K94").Value = 22728444
$ws.Range("L94").Value = 1955
$ws.Range("M94").Value = -22727993
$ws.Range("N94").Value = -2857
$ws.Range("H107").Value = 1676.75
$ws.Range("I107").Value = 1395.4706
$ws.Range("K107").Value = 1395.4706
$ws.Range("M107").Value = 524.5293999999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2623.889
$ws.Range("I31").Value = 3354.9092
$ws.Range("K31").Value = 3354.9092
$ws.Range("M31").Value = -3059.9092
$ws.Range("H34").Value = 2623.889
$ws.Range("I34").Value = 3354.9092
$ws.Range("K34").Value = 3354.9092
$ws.Range("M34").Value = -3152.9092
$ws.Range("H109").Value = 16612.625
$ws.Range("J109").Value = 16612.625
$ws.Range("L109").Value = 16612.625
$ws.Range("N109").Value = -18692.625
$ws.Range("H132").Value = 2858.5908
$ws.Range("I132").Value = 2568.75
$ws.Range("J132").Value = 5757
$ws.Range("K132").Value = 7706.25
$ws.Range("L132").Value = 17271
$ws.Range("M132").Value = -5176.25
$ws.Range("N132").Value = -22331

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2683.9167
$ws.Range("I5").Value = 3636.7144
$ws.Range("J5").Value = 1350
$ws.Range("K5").Value = 10910.1432
$ws.Range("L5").Value = 4050
$ws.Range("M5").Value = -10798.1432
$ws.Range("N5").Value = -4274
$ws.Range("H11").Value = 147682.42
$ws.Range("I11").Value = 172264.5
$ws.Range("J11").Value = 190
$ws.Range("K11").Value = 516793.5
$ws.Range("L11").Value = 570
$ws.Range("M11").Value = -516653.5
$ws.Range("N11").Value = -850
$ws.Range("H14").Value = 144.7619
$ws.Range("I14").Value = 144.7619
$ws.Range("K14").Value = 434.2857
$ws.Range("M14").Value = -261.2857
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3000
$ws.Range("L22").ClearContents()
$ws.Range("M22").Value = -2831
$ws.Range("N22").Value = 0
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 3000
$ws.Range("L27").ClearContents()
$ws.Range("M27").Value = -2898
$ws.Range("N27").Value = 0
$ws.Range("H29").Value = 575.5
$ws.Range("I29").Value = 100
$ws.Range("J29").Value = 734
$ws.Range("K29").Value = 300
$ws.Range("L29").Value = 2202
$ws.Range("M29").Value = -23
$ws.Range("N29").Value = -2756
$ws.Range("H107").Value = 10551.9
$ws.Range("I107").Value = 606.6667
$ws.Range("J107").Value = 14814.143
$ws.Range("K107").Value = 1820.0001
$ws.Range("L107").Value = 44442.429
$ws.Range("M107").Value = 99.99990000000003
$ws.Range("N107").Value = -48282.429
$ws.Range("H122").Value = 794.64
$ws.Range("I122").Value = 617
$ws.Range("J122").Value = 934.2143
$ws.Range("K122").Value = 5553
$ws.Range("L122").Value = 8407.9287
$ws.Range("M122").Value = -3103
$ws.Range("N122").Value = -13307.9287
$ws.Range("H131").Value = 15626200
$ws.Range("I131").Value = 200000400
$ws.Range("J131").Value = 1267.5593
$ws.Range("K131").Value = 600001200
$ws.Range("L131").Value = 3802.6779
$ws.Range("M131").Value = -599996160
$ws.Range("N131").Value = -13882.6779
$ws.Range("H135").Value = 2683.9167
$ws.Range("I135").Value = 3636.7144
$ws.Range("J135").Value = 1350
$ws.Range("K135").Value = 32730.4296
$ws.Range("L135").Value = 12150
$ws.Range("M135").Value = -30195.4296
$ws.Range("N135").Value = -17220
$ws.Range("H137").Value = 20838580
$ws.Range("I137").Value = 46876492
$ws.Range("J137").Value = 8248.85
$ws.Range("K137").Value = 140629476
$ws.Range("L137").Value = 24746.55
$ws.Range("M137").Value = -140624376
$ws.Range("N137").Value = -34946.55

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").ClearContents()
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = 0
$ws.Range("H97").Value = 454.08334
$ws.Range("I97").Value = 468.0909
$ws.Range("K97").Value = 468.0909
$ws.Range("M97").Value = 27.90910000000002
$ws.Range("H99").Value = 4688.3335
$ws.Range("I99").Value = 1626
$ws.Range("J99").Value = 20000
$ws.Range("K99").Value = 1626
$ws.Range("L99").Value = 20000
$ws.Range("M99").Value = 620
$ws.Range("N99").Value = -24492
$ws.Range("H113").Value = 1136.4286
$ws.Range("I113").Value = 1211
$ws.Range("J113").Value = 950
$ws.Range("K113").Value = 1211
$ws.Range("L113").Value = 950
$ws.Range("M113").Value = 959
$ws.Range("N113").Value = -5290
$ws.Range("H126").Value = 2210.1875
$ws.Range("I126").Value = 1831.5
$ws.Range("J126").Value = 2437.4
$ws.Range("K126").Value = 5494.5
$ws.Range("L126").Value = 7312.200000000001
$ws.Range("M126").Value = -3024.5
$ws.Range("N126").Value = -12252.2
$ws.Range("H132").Value = 1623.7142
$ws.Range("I132").Value = 1458.6
$ws.Range("K132").Value = 4375.799999999999
$ws.Range("M132").Value = -1845.799999999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1899.6666
$ws.Range("I61").Value = 1649.5
$ws.Range("J61").Value = 2400
$ws.Range("K61").Value = 1649.5
$ws.Range("L61").Value = 2400
$ws.Range("M61").Value = -1447.5
$ws.Range("N61").Value = -2804
$ws.Range("H113").Value = 1899.6666
$ws.Range("I113").Value = 1649.5
$ws.Range("J113").Value = 2400
$ws.Range("K113").Value = 1649.5
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = 520.5
$ws.Range("N113").Value = -6740

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 310.8889
$ws.Range("I81").Value = 324.75
$ws.Range("J81").Value = 200
$ws.Range("K81").Value = 649.5
$ws.Range("L81").Value = 400
$ws.Range("M81").Value = 411.5
$ws.Range("N81").Value = -2522
$ws.Range("H84").Value = 310.8889
$ws.Range("I84").Value = 324.75
$ws.Range("J84").Value = 200
$ws.Range("K84").Value = 3247.5
$ws.Range("L84").Value = 2000
$ws.Range("M84").Value = 2056.5
$ws.Range("N84").Value = -12608
$ws.Range("H109").Value = 35085.5
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774
$ws.Range("H132").Value = 7176.273
$ws.Range("I132").Value = 7643.9
$ws.Range("K132").Value = 22931.7
$ws.Range("M132").Value = -20401.7
$ws.Range("H136").Value = 609.7368
$ws.Range("I136").Value = 434.96667
$ws.Range("J136").Value = 1265.125
$ws.Range("K136").Value = 1304.90001
$ws.Range("L136").Value = 3795.375
$ws.Range("M136").Value = 1245.09999
$ws.Range("N136").Value = -8895.375
